# LR3/table_1_109.xlsx - "Updated file (5) 11.10.22"
#
# Semantic changes applied:
#  - D4 formula changes from "=$A$1*1.1" to "=D3"
#  - D5:D34 become a copy-down of the row above ("=D4", "=D5", ... "=D33")
#  - H4 formula changes from "=IF(G3+1<F3+1,H3,H3+1)" to "=IF(G3>F3,G3-F3,0)"
#  - H5:H38 become a copy-down of the same new pattern based on the row above
#  - C40 formula changes from "=TRUNC(SUM(K3:K38))" to "=FLOOR(SUM(K3:K38),1)"
#  - Selection moves to D5 (mirrors the saved cursor position in the file)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column D: "price" recompute switched from the absolute A1 formula to a
#     simple copy-down of the previous row ---
$ws.Range("D4").Formula = "=D3"

for ($r = 5; $r -le 34; $r++) {
    $prev = $r - 1
    $ws.Range("D$r").Formula = "=D$prev"
}

# --- Column H: overdue-days logic rewritten ---
$ws.Range("H4").Formula = "=IF(G3>F3,G3-F3,0)"

for ($r = 5; $r -le 38; $r++) {
    $prev = $r - 1
    $ws.Range("H$r").Formula = "=IF(G$prev>F$prev,G$prev-F$prev,0)"
}

# --- Totals: TRUNC -> FLOOR ---
$ws.Range("C40").Formula = "=FLOOR(SUM(K3:K38),1)"

# --- Cursor / selection bookkeeping (matches the saved file's selection) ---
$ws.Range("D5").Select()
